$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of book data to append (Title, Author, Start Date serial, Finish Date serial, Tags, Type, Length)
$rows = @(
    @("Infinite Powers", "Steven Strogatz", 43881, 43897, "math;calculus;infinity;history;newton", "Hard Copy", "301 Pages"),
    @("The 4-Hour Work Week", "Tim Ferris", 43895, 43897, "fullfillment;self improvement;productivity", "Audio", "13 Hours 1 Min"),
    @("Elizabeth II: Life of a Monarch", "Ruth Cowen", 43901, 43902, "queen Elizabeth;monarchy;biography;england;history", "Audio", "3 Hours 47 Mins")
)

$lastExistingDateRow = 34
$startRow = 35

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]

    # Set raw date serials first (avoids the engine minting a fresh custom
    # number format when a date-like string is assigned), then copy the
    # date number formatting from an existing date cell so the new cells
    # share the workbook's existing date style instead of a new one.
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]

    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
}

$dateRange = $ws.Range($ws.Cells.Item($lastExistingDateRow, 3), $ws.Cells.Item($lastExistingDateRow, 4))
$dateRange.Copy()
$newDateRange = $ws.Range($ws.Cells.Item($startRow, 3), $ws.Cells.Item($startRow + $rows.Length - 1, 4))
$newDateRange.PasteSpecial(-4122)
$excel.CutCopyMode = 0

$lastRow = $startRow + $rows.Length - 1
$ws.Range("A" + ($lastRow + 1)).Select()
$excel.ActiveWindow.ScrollRow = 18
